# Updated capital structure database
# - Refreshes the financial metrics for the existing Brazil "Banks (Regional)"
#   rows (2-4), renames row 3 to "Banco do Estado de Sergipe S.A." and bumps
#   the company id in row 2 from "2" to "4".
# - Appends two additional companies as new rows 5 and 6 ("Banestes S.A -
#   Banco do Estado do Espirito Santo", previously in row 3, and "Banco Alfa
#   de Investimento S.A."), extending the sheet's used range to A1:AQ6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 2: update existing entity (company id "2" -> "4") and metrics ----
# Row 2
$ws.Range("A2").Value = "Brazil"
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "4"
$ws.Range("B2").ClearFormats()
$ws.Range("C2").Value = "Banks (Regional)"
$ws.Range("D2").Value = 0.08635000000000001
$ws.Range("E2").Value = 0.0418
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 122.8
$ws.Range("L2").Value = 0.2073273678878947
$ws.Range("M2").Value = 46.5245
$ws.Range("N2").Value = 0.05549862817607062
$ws.Range("O2").Value = 0.378864006514658
$ws.Range("P2").Value = 46.4785
$ws.Range("Q2").Value = 0.05544375521889539
$ws.Range("R2").Value = 0.3784894136807818
$ws.Range("S2").Value = 0.04599999999999982
$ws.Range("T2").Value = 0.000988726369977105
$ws.Range("U2").Value = 1238.4
$ws.Range("V2").Value = 1.477275438387212
$ws.Range("W2").Value = 0.1102233645011941
$ws.Range("X2").Value = 0.1898079483513983
$ws.Range("Y2").Value = -0.07958458385020419
$ws.Range("Z2").Value = 0.09907333065703197
$ws.Range("AA2").Value = 0
$ws.Range("AB2").Value = 0.05566985909404883
$ws.Range("AC2").Value = -0.05566985909404883
$ws.Range("AD2").Value = 4875.8
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 4875.8
$ws.Range("AG2").Value = 3637.4
$ws.Range("AH2").Value = 0.8532927320137904
$ws.Range("AI2").Value = 0.819613710097665
$ws.Range("AJ2").Value = 0.812699689434055
$ws.Range("AK2").Value = 0.7721897887697697
$ws.Range("AL2").Value = 0
$ws.Range("AM2").Value = 0

# Row 3
$ws.Range("A3").Value = "Brazil"
$ws.Range("B3").Value = "Banco do Estado de Sergipe S.A. (BOVESPA:BGIP4)"
$ws.Range("C3").Value = "Banks (Regional)"
$ws.Range("D3").Value = 0.233
$ws.Range("E3").Value = 0.188
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 10.8
$ws.Range("L3").Value = 0.1003717472118959
$ws.Range("M3").Value = 3.7485
$ws.Range("N3").Value = 0.03614754098360656
$ws.Range("O3").Value = 0.3470833333333333
$ws.Range("P3").Value = 3.7485
$ws.Range("Q3").Value = 0.03614754098360656
$ws.Range("R3").Value = 0.3470833333333333
$ws.Range("S3").Value = 0
$ws.Range("T3").Value = 0
$ws.Range("U3").Value = 237.3
$ws.Range("V3").Value = 2.288331726133076
$ws.Range("W3").Value = 0.1181619256017506
$ws.Range("X3").Value = 0.06077370760841536
$ws.Range("Y3").Value = 0.0573882179933352
$ws.Range("Z3").Value = -1.761047463175124
$ws.Range("AA3").Value = -0
$ws.Range("AB3").Value = 0.05247872550769265
$ws.Range("AC3").Value = -0.05247872550769265
$ws.Range("AD3").Value = 38.6
$ws.Range("AE3").Value = 0
$ws.Range("AF3").Value = 38.6
$ws.Range("AG3").Value = -198.7
$ws.Range("AH3").Value = 0.2712579058327477
$ws.Range("AI3").Value = 0.2852919438285292
$ws.Range("AJ3").Value = 2.091578947368421
$ws.Range("AK3").Value = 1.948039215686274
$ws.Range("AL3").Value = 0
$ws.Range("AM3").Value = 0

# Row 4
$ws.Range("A4").Value = "Brazil"
$ws.Range("B4").Value = "Banco da Amazônia S.A. (BOVESPA:BAZA3)"
$ws.Range("C4").Value = "Banks (Regional)"
$ws.Range("D4").Value = -0.0122
$ws.Range("E4").Value = 0.0132
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 60.3
$ws.Range("L4").Value = 0.279814385150812
$ws.Range("M4").Value = 18.9
$ws.Range("N4").Value = 0.0784557907845579
$ws.Range("O4").Value = 0.3134328358208955
$ws.Range("P4").Value = 18.9
$ws.Range("Q4").Value = 0.0784557907845579
$ws.Range("R4").Value = 0.3134328358208955
$ws.Range("S4").Value = 0
$ws.Range("T4").Value = 0
$ws.Range("U4").Value = 183.8
$ws.Range("V4").Value = 0.7629721876297219
$ws.Range("W4").Value = 0.1230361150785554
$ws.Range("X4").Value = 0.08775538275785405
$ws.Range("Y4").Value = 0.03528073232070134
$ws.Range("Z4").Value = 0.6061884669479605
$ws.Range("AA4").Value = 0
$ws.Range("AB4").Value = 0.05456830755983145
$ws.Range("AC4").Value = -0.05456830755983145
$ws.Range("AD4").Value = 328
$ws.Range("AE4").Value = 0
$ws.Range("AF4").Value = 328
$ws.Range("AG4").Value = 144.2
$ws.Range("AH4").Value = 0.5765512392336087
$ws.Range("AI4").Value = 0.4449871116537783
$ws.Range("AJ4").Value = 0.3744481952739548
$ws.Range("AK4").Value = 0.2606181095246702
$ws.Range("AL4").Value = 0
$ws.Range("AM4").Value = 0

# Row 5
$ws.Range("A5").Value = "Brazil"
$ws.Range("B5").Value = "Banestes S.A - Banco do Estado do Espírito Santo (BOVESPA:BEES3)"
$ws.Range("C5").Value = "Banks (Regional)"
$ws.Range("D5").Value = 0.0687
$ws.Range("E5").Value = 0.0704
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 38.5
$ws.Range("L5").Value = 0.1765245300320954
$ws.Range("M5").Value = 20.4
$ws.Range("N5").Value = 0.05982404692082111
$ws.Range("O5").Value = 0.5298701298701298
$ws.Range("P5").Value = 20.4
$ws.Range("Q5").Value = 0.05982404692082111
$ws.Range("R5").Value = 0.5298701298701298
$ws.Range("S5").Value = 0
$ws.Range("T5").Value = 0
$ws.Range("U5").Value = 61.3
$ws.Range("V5").Value = 0.1797653958944281
$ws.Range("W5").Value = 0.1022848034006376
$ws.Range("X5").Value = 0.2918605139449425
$ws.Range("Y5").Value = -0.1895757105443049
$ws.Range("Z5").Value = 0.0852019689038206
$ws.Range("AA5").Value = 0
$ws.Range("AB5").Value = 0.0567714106282662
$ws.Range("AC5").Value = -0.0567714106282662
$ws.Range("AD5").Value = 3016.3
$ws.Range("AE5").Value = 0
$ws.Range("AF5").Value = 3016.3
$ws.Range("AG5").Value = 2955
$ws.Range("AH5").Value = 0.8984302862419206
$ws.Range("AI5").Value = 0.9118198307134221
$ws.Range("AJ5").Value = 0.8965412621359223
$ws.Range("AK5").Value = 0.9101549265407953
$ws.Range("AL5").Value = 0
$ws.Range("AM5").Value = 0

# Row 6
$ws.Range("A6").Value = "Brazil"
$ws.Range("B6").Value = "Banco Alfa de Investimento S.A. (BOVESPA:BRIV4)"
$ws.Range("C6").Value = "Banks (Regional)"
$ws.Range("D6").Value = 0.104
$ws.Range("E6").Value = -0.008659999999999999
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 13.2
$ws.Range("L6").Value = 0.2583170254403131
$ws.Range("M6").Value = 3.476
$ws.Range("N6").Value = 0.02276358873608383
$ws.Range("O6").Value = 0.2633333333333334
$ws.Range("P6").Value = 3.43
$ws.Range("Q6").Value = 0.02246234446627374
$ws.Range("R6").Value = 0.2598484848484849
$ws.Range("S6").Value = 0.04599999999999982
$ws.Range("T6").Value = 0.01323360184119673
$ws.Range("U6").Value = 756
$ws.Range("V6").Value = 4.950884086444008
$ws.Range("W6").Value = 0.03450078410872974
$ws.Range("X6").Value = 0.317257611885869
$ws.Range("Y6").Value = -0.2827568277771392
$ws.Range("Z6").Value = 0.01635618718391909
$ws.Range("AA6").Value = 0
$ws.Range("AB6").Value = 0.05683148356705298
$ws.Range("AC6").Value = -0.05683148356705298
$ws.Range("AD6").Value = 1492.9
$ws.Range("AE6").Value = 0
$ws.Range("AF6").Value = 1492.9
$ws.Range("AG6").Value = 736.9000000000001
$ws.Range("AH6").Value = 0.9072070977151191
$ws.Range("AI6").Value = 0.8441617189708793
$ws.Range("AJ6").Value = 0.8283498201438849
$ws.Range("AK6").Value = 0.7278024691358025
$ws.Range("AL6").Value = 0
$ws.Range("AM6").Value = 0

Write-Output "Capital structure database updated: rows 2-6 refreshed, 2 new rows added"
